# Update "gh-pages" scraped output at 456a3b4
# Sheet 1 = "展览" (index 1), Sheet 4 = "全部类型" (index 4) both carry the
# full combined-event table and receive (mostly) the same F-column
# ("想去人数" / want-to-go count) refresh plus one updated cover image URL.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet "展览" (index 1) updates ---
$wsExhibit.Range("F2").Value  = 285
$wsExhibit.Range("F3").Value  = 1453
$wsExhibit.Range("F4").Value  = 178
$wsExhibit.Range("F6").Value  = 241
$wsExhibit.Range("F10").Value = 144
$wsExhibit.Range("F11").Value = 6
$wsExhibit.Range("F12").Value = 4817
$wsExhibit.Range("F14").Value = 7094
$wsExhibit.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202407/ptM44zuA1721386766857.jpeg"
$wsExhibit.Range("F20").Value = 13
$wsExhibit.Range("F21").Value = 4201
$wsExhibit.Range("F22").Value = 1394
$wsExhibit.Range("F23").Value = 86
$wsExhibit.Range("F24").Value = 80
$wsExhibit.Range("F25").Value = 2775
$wsExhibit.Range("F28").Value = 181
$wsExhibit.Range("F29").Value = 401
$wsExhibit.Range("F31").Value = 420
$wsExhibit.Range("F32").Value = 253
$wsExhibit.Range("F34").Value = 1659
$wsExhibit.Range("F35").Value = 1082
$wsExhibit.Range("F36").Value = 77
$wsExhibit.Range("F37").Value = 993
$wsExhibit.Range("F38").Value = 92
$wsExhibit.Range("F41").Value = 503
$wsExhibit.Range("F42").Value = 12
$wsExhibit.Range("F43").Value = 29
$wsExhibit.Range("F45").Value = 1184
$wsExhibit.Range("F46").Value = 662

# --- Sheet "全部类型" (index 4) updates ---
$wsAll.Range("F2").Value  = 285
$wsAll.Range("F3").Value  = 1453
$wsAll.Range("F4").Value  = 178
$wsAll.Range("F6").Value  = 241
$wsAll.Range("F10").Value = 144
$wsAll.Range("F11").Value = 6
$wsAll.Range("F12").Value = 4817
$wsAll.Range("F14").Value = 7094
$wsAll.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202407/ptM44zuA1721386766857.jpeg"
$wsAll.Range("F18").Value = 586
$wsAll.Range("F20").Value = 13
$wsAll.Range("F21").Value = 4201
$wsAll.Range("F22").Value = 1394
$wsAll.Range("F23").Value = 86
$wsAll.Range("F24").Value = 80
$wsAll.Range("F25").Value = 2775
$wsAll.Range("F28").Value = 181
$wsAll.Range("F29").Value = 401
$wsAll.Range("F31").Value = 420
$wsAll.Range("F32").Value = 253
$wsAll.Range("F34").Value = 1659
$wsAll.Range("F35").Value = 1082
$wsAll.Range("F36").Value = 77
$wsAll.Range("F37").Value = 993
$wsAll.Range("F38").Value = 92
$wsAll.Range("F41").Value = 503
$wsAll.Range("F42").Value = 12
$wsAll.Range("F43").Value = 29
$wsAll.Range("F45").Value = 1184
$wsAll.Range("F46").Value = 662
